$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.155.42"
$ws.Range("E2").Value = "  +1.61%  "
$ws.Range("D3").Value = "3.381.79"
$ws.Range("E3").Value = "  +1.16%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'586.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("D6").Value = "'179.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.98%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  +0.75%  "
$ws.Range("D9").Value = "'0.193"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.89%  "
$ws.Range("D10").Value = "'0.590"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.21%  "
$ws.Range("D11").Value = "'48.38"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.28%  "
$ws.Range("D12").Value = "'0.0000280"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.92%  "
$ws.Range("D13").Value = "'679.05"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.10%  "
$ws.Range("D14").Value = "'8.60"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.62%  "
$ws.Range("D15").Value = "3.927.97"
$ws.Range("E15").Value = "  +1.20%  "
$ws.Range("D16").Value = "69.223.15"
$ws.Range("E16").Value = "  +1.64%  "
$ws.Range("E17").Value = "  +1.65%  "
$ws.Range("D18").Value = "3.375.26"
$ws.Range("E18").Value = "  +1.00%  "
$ws.Range("D19").Value = "'17.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.89%  "
$ws.Range("D20").Value = "'11.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.65%  "
$ws.Range("D21").Value = "'0.900"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").Value = "'5.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.50%  "
$ws.Range("D23").Value = "'17.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.30%  "
$ws.Range("D24").Value = "'102.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.67%  "
$ws.Range("D25").Value = "'3.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("D26").Value = "'2.71"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.16%  "
$ws.Range("D27").Value = "'9.57"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.17%  "
$ws.Range("D28").Value = "'34.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.91%  "
$ws.Range("D29").Value = "'8.69"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.03%  "
$ws.Range("D30").Value = "'6.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.42%  "
$ws.Range("D31").Value = "'11.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.20%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "'553.48"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.18%  "
$ws.Range("B33").Value = "dogwifhat"
$ws.Range("C33").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D33").Value = "'3.61"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.25%  "
$ws.Range("E34").Value = "  +0.22%  "
$ws.Range("D35").Value = "'57.97"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.17%  "
$ws.Range("D37").Value = "3.680.09"
$ws.Range("E37").Value = "  -1.17%  "
$ws.Range("D38").Value = "'0.139"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.74%  "
$ws.Range("D39").Value = "'34.96"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.57%  "
$ws.Range("D40").Value = "'3.23"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.26%  "
$ws.Range("D41").Value = "'2.66"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.46%  "
$ws.Range("D42").Value = "0.0₃0696"
$ws.Range("E42").Value = "  +2.49%  "
$ws.Range("D43").Value = "'0.337"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.15%  "
$ws.Range("D44").Value = "'0.0422"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.48%  "
$ws.Range("E45").Value = "  -1.93%  "
$ws.Range("E46").Value = "  -0.68%  "
$ws.Range("E47").Value = "  +0.53%  "
$ws.Range("E48").Value = "  +5.18%  "
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("D50").Value = "'131.57"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.25%  "
$ws.Range("D51").Value = "'2.57"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.24%  "
